# Add the "Added checkout information test case" worksheet (UserInfo)
# after the existing "login" sheet, and populate it with sample
# checkout / user-info data used by the test fixture.

$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item(1)

# New sheet goes right after "login"
$ws = $wb.Worksheets.Add($null, $loginSheet)
$ws.Name = "UserInfo"

# Header row
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "ZIPCode"

# Data rows
$ws.Range("A2").Value = "abc"
$ws.Range("B2").Value = "xyz"
$ws.Range("C2").Value = 9999

$ws.Range("A3").Value = "ABC"
$ws.Range("B3").Value = "XYZ"
$ws.Range("C3").Value = 8888

# Auto-fit column A like Excel does when a user double-clicks the
# column border after entering data.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the selection where the author last clicked, and make the
# new sheet the active / visible tab (login loses tabSelected).
$ws.Range("G14").Select() | Out-Null
$ws.Activate() | Out-Null
